# Apply the two changes recorded in the commit:
#  1. The table on slide 6 gets a new built-in table style applied.
#  2. The deck's theme colour scheme is changed from the "Integral" palette
#     to the standard Office palette (the font/format schemes were already
#     identical between the two themes, so only the colour scheme differs).

$p = $ppt.ActivePresentation

# --- 1. Change the table style on slide 6 -----------------------------
$s = $p.Slides.Item(6)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{AFF17EAD-9C35-4E8E-AA47-6333B77AFAF0}")
    }
}

# --- 2. Re-colour the theme from "Integral" to the default "Office" ----
function RGBVal($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

# order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeColors = @(
    @(0x00, 0x00, 0x00),  # dk1
    @(0xFF, 0xFF, 0xFF),  # lt1
    @(0x44, 0x54, 0x6A),  # dk2
    @(0xE7, 0xE6, 0xE6),  # lt2
    @(0x5B, 0x9B, 0xD5),  # accent1
    @(0xED, 0x7D, 0x31),  # accent2
    @(0xA5, 0xA5, 0xA5),  # accent3
    @(0xFF, 0xC0, 0x00),  # accent4
    @(0x44, 0x72, 0xC4),  # accent5
    @(0x70, 0xAD, 0x47),  # accent6
    @(0x05, 0x63, 0xC1),  # hlink
    @(0x95, 0x4F, 0x72)   # folHlink
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $rgb = $officeColors[$i - 1]
    $tcs.Colors($i).RGB = RGBVal $rgb[0] $rgb[1] $rgb[2]
}
